# Commit: "Added column for CNN benchmark"
# - rename the "CNN" header (J1) to "CNN logloss"
# - populate the new CNN logloss values in J2:J21 (previously empty, which
#   made the AVERAGE(J2:J21) in J22 blow up with #DIV/0!)
# - give column J a wider custom width so the longer header fits
# - leave the sheet-view selection where the author last left it (J24)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text
$ws.Range("J1").Value = "CNN logloss"

# New CNN logloss benchmark values for each row
$ws.Range("J2").Value = 5.7
$ws.Range("J3").Value = 17.76
$ws.Range("J4").Value = 18.3
$ws.Range("J5").Value = 4.79
$ws.Range("J6").Value = 1.55
$ws.Range("J7").Value = 3.64
$ws.Range("J8").Value = 8.42
$ws.Range("J9").Value = 15.07
$ws.Range("J10").Value = 11.97
$ws.Range("J11").Value = 5.63
$ws.Range("J12").Value = 4.79
$ws.Range("J13").Value = 9.88
$ws.Range("J14").Value = 4.59
$ws.Range("J15").Value = 5.62
$ws.Range("J16").Value = 13.63
$ws.Range("J17").Value = 0.829
$ws.Range("J18").Value = 15.69
$ws.Range("J19").Value = 0.34
$ws.Range("J20").Value = 2.97
$ws.Range("J21").Value = 15.42

# The J22 average formula already exists (AVERAGE(J2:J21)); now that J2:J21
# are populated it naturally recalculates away from the #DIV/0! error.

# Widen column J to fit the new header/values (target OOXML width ~21.5
# characters; ColumnWidth is stored in pixel-quantized character units, so
# feed it the value whose quantized result lands closest to 21.5)
$ws.Columns.Item(10).ColumnWidth = 20.71

# Restore the author's last selected cell
$ws.Range("J24").Select()
